# Update the "Förändrad" (Changed) date column (C) for rows 2-5
# from serial date 45185 (2023-09-16) to serial date 45204 (2023-10-05).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C5").Value2 = 45204
